$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 532, pushing the existing rows 532:644 down to 534:646.
$ws.Rows("532:533").Insert()

# Populate new row 532 (Packham's Triumph / Primera)
$ws.Range("A532").Value = 4
$ws.Range("B532").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C532").Value = "Los Lagos"
$ws.Range("D532").Value = 45275
$ws.Range("D532").NumberFormat = $ws.Range("D534").NumberFormat
$ws.Range("E532").Value = 10
$ws.Range("F532").Value = "Fruta"
$ws.Range("G532").Value = 100104
$ws.Range("H532").Value = "Frutos de pepita"
$ws.Range("I532").Value = 100104005
$ws.Range("J532").Value = "Pera"
$ws.Range("K532").Value = "Packham's Triumph"
$ws.Range("L532").Value = "Primera"
$ws.Range("M532").Value = 300
$ws.Range("N532").Value = 21000
$ws.Range("O532").Value = 21000
$ws.Range("P532").Value = 21000
$ws.Range("Q532").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R532").Value = "Región de O'Higgins"
$ws.Range("S532").Value = 1400
$ws.Range("T532").Value = 15

# Populate new row 533 (Packham's Triumph / Segunda)
$ws.Range("A533").Value = 4
$ws.Range("B533").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C533").Value = "Los Lagos"
$ws.Range("D533").Value = 45275
$ws.Range("D533").NumberFormat = $ws.Range("D534").NumberFormat
$ws.Range("E533").Value = 10
$ws.Range("F533").Value = "Fruta"
$ws.Range("G533").Value = 100104
$ws.Range("H533").Value = "Frutos de pepita"
$ws.Range("I533").Value = 100104005
$ws.Range("J533").Value = "Pera"
$ws.Range("K533").Value = "Packham's Triumph"
$ws.Range("L533").Value = "Segunda"
$ws.Range("M533").Value = 300
$ws.Range("N533").Value = 17000
$ws.Range("O533").Value = 17000
$ws.Range("P533").Value = 17000
$ws.Range("Q533").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R533").Value = "Región de O'Higgins"
$ws.Range("S533").Value = 1133
$ws.Range("T533").Value = 15
